# Adds a method to retrieve the fundamentals from the stock: appends a new
# trade record (row 15) to the CELG trade history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 9531.77
$ws.Range("B15").Value = 9339.3799999999992
$ws.Range("C15").Value = 104.49
$ws.Range("D15").Value = 106.64
$ws.Range("E15").Value = $false
$ws.Range("F15").Value = 2.06
$ws.Range("G15").Value = 42626.54446759259
$ws.Range("H15").Value = $true
